# ToDo.xlsx update
# - Mark "Statistik" (row 9) and "Benutzer bearbeiten" (row 10) as done (Erledigt? = X),
#   which makes the existing AutoFilter (hides "done" rows - filters on blank) hide them.
# - Mark row 13 ("Alle Warnings ...") as done too, and note a new follow-up task in the
#   "Problem" column.
# - Rename the "Datenbank online bringen" task to the more complete description.
# - Grow the used range / AutoFilter / _FilterDatabase by one row (A1:F16 -> A1:F17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Statistik: mark as erledigt and hide (filtered out)
$ws.Range("D9").Value = "X"
$ws.Rows.Item(9).Hidden = $true

# Row 10 - Benutzer bearbeiten: mark as erledigt and hide (filtered out)
$ws.Range("D10").Value = "X"
$ws.Rows.Item(10).Hidden = $true

# Row 14 - rename task text (not done yet, stays visible).
$ws.Range("B14").Value = "Datenbank und Rest online bringen (multiplechoice.szüsz.de)"

# Row 13 - Alle Warnings und Fehler nicht auf der Seite anzeigen: mark as erledigt, hide,
# and add a follow-up note in the "Problem" column.
$ws.Range("D13").Value = "X"
$ws.Range("E13").Value = "Suchen nach zu 0"
$ws.Rows.Item(13).Hidden = $true

# Grow the filtered/used range from A1:F16 to A1:F17 and re-apply the
# "blank" filter on the "Erledigt?" column (column D, colId=3).
$ws.AutoFilterMode = $false
$ws.Range("A1:F17").AutoFilter(4, @(""), 7)

# Keep the hidden _FilterDatabase defined name pointing at the resized range.
$wb.Names.Item(1).RefersTo = "=ToDo!`$A`$1:`$F`$17"
